# Fruta / hortaliza, semanal
# Insert two new weekly records (Santina, week of 2021-11-18) into the
# "Vega Monumental Concepción - Cereza" dataset right before the existing
# row 38 (Bing), shifting the rest of the Cereza block down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 38..51 down to 40..53, leaving two blank rows (38:39) for the
# new "Santina" entries. Excel inherits the row-38 formatting (incl. the
# date style on column D) for the freshly inserted rows.
$ws.Rows("38:39").Insert()

# --- New row 38: Santina / Primera-equivalent record ---
$ws.Range("A38").Value2 = 11
$ws.Range("B38").Value2 = "Vega Monumental Concepción"
$ws.Range("C38").Value2 = "Bíobío"
$ws.Range("D38").Value2 = 44518
$ws.Range("E38").Value2 = 8
$ws.Range("F38").Value2 = "Fruta"
$ws.Range("G38").Value2 = 100103
$ws.Range("H38").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I38").Value2 = 100103001
$ws.Range("J38").Value2 = "Cereza"
$ws.Range("K38").Value2 = "Santina"
$ws.Range("L38").Value2 = "Primera"
$ws.Range("M38").Value2 = 50
$ws.Range("N38").Value2 = 32000
$ws.Range("O38").Value2 = 32000
$ws.Range("P38").Value2 = 32000
$ws.Range("Q38").Value2 = "$/caja 10 kilos"
$ws.Range("R38").Value2 = "Provincia de Curicó"
$ws.Range("S38").Value2 = 3200
$ws.Range("T38").Value2 = 10

# --- New row 39: Santina / Segunda-equivalent record ---
$ws.Range("A39").Value2 = 11
$ws.Range("B39").Value2 = "Vega Monumental Concepción"
$ws.Range("C39").Value2 = "Bíobío"
$ws.Range("D39").Value2 = 44518
$ws.Range("E39").Value2 = 8
$ws.Range("F39").Value2 = "Fruta"
$ws.Range("G39").Value2 = 100103
$ws.Range("H39").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I39").Value2 = 100103001
$ws.Range("J39").Value2 = "Cereza"
$ws.Range("K39").Value2 = "Santina"
$ws.Range("L39").Value2 = "Segunda"
$ws.Range("M39").Value2 = 50
$ws.Range("N39").Value2 = 28000
$ws.Range("O39").Value2 = 28000
$ws.Range("P39").Value2 = 28000
$ws.Range("Q39").Value2 = "$/caja 10 kilos"
$ws.Range("R39").Value2 = "Provincia de Curicó"
$ws.Range("S39").Value2 = 2800
$ws.Range("T39").Value2 = 10
